# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") holds a per-row statistic (formerly a strike count, now the
# freshly calculated/simulated "K" value). The workbook stores these as plain
# literal numbers (no formula), so we rewrite the newly computed values cell
# by cell on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 2
    5  = 4
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 0
    13 = 1
    14 = 3
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 2
    22 = 2
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 2
    30 = 3
    31 = 1
    32 = 1
    33 = 2
    34 = 1
    35 = 0
    36 = 1
    37 = 1
    38 = 2
    39 = 1
    40 = 2
    41 = 2
    42 = 1
    44 = 0
    45 = 2
    46 = 2
    47 = 1
    48 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}

Write-Output "updated $($newK.Count) K values in column G"
